$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (155) down to each new row,
# then fill in the values. PasteSpecial with xlPasteFormats (-4122) copies only
# formats so it reuses the existing cell style (s="2" on col A) instead of minting a new one.
$ws.Range("A155:G155").Copy()

$ws.Range("A156:G156").PasteSpecial(-4122)
$ws.Range("A156").Value = 44113
$ws.Range("B156").Value = "10:00:00"
$ws.Range("C156").Value = 3027
$ws.Range("D156").Value = 1515
$ws.Range("E156").Value = 109
$ws.Range("F156").Value = 2664
$ws.Range("G156").Value = 254

$ws.Range("A157:G157").PasteSpecial(-4122)
$ws.Range("A157").Value = 44116
$ws.Range("B157").Value = "11:15:00"
$ws.Range("C157").Value = 3167
$ws.Range("D157").Value = 1585
$ws.Range("E157").Value = 109
$ws.Range("F157").Value = 2723
$ws.Range("G157").Value = 335

$ws.Range("A158:G158").PasteSpecial(-4122)
$ws.Range("A158").Value = 44117
$ws.Range("B158").Value = "11:45:00"
$ws.Range("C158").Value = 3188
$ws.Range("D158").Value = 1594
$ws.Range("E158").Value = 109
$ws.Range("F158").Value = 2772
$ws.Range("G158").Value = 307

$ws.Range("A159:G159").PasteSpecial(-4122)
$ws.Range("A159").Value = 44118
$ws.Range("B159").Value = "11:45:00"
$ws.Range("C159").Value = 3257
$ws.Range("D159").Value = 1626
$ws.Range("E159").Value = 109
$ws.Range("F159").Value = 2806
$ws.Range("G159").Value = 342

$ws.Range("A160:G160").PasteSpecial(-4122)
$ws.Range("A160").Value = 44119
$ws.Range("B160").Value = "10:15:00"
$ws.Range("C160").Value = 3393
$ws.Range("D160").Value = 1672
$ws.Range("E160").Value = 109
$ws.Range("F160").Value = 2845
$ws.Range("G160").Value = 439

$ws.Range("A161:G161").PasteSpecial(-4122)
$ws.Range("A161").Value = 44120
$ws.Range("B161").Value = "12:15:00"
$ws.Range("C161").Value = 3509
$ws.Range("D161").Value = 1729
$ws.Range("E161").Value = 109
$ws.Range("F161").Value = 2891
$ws.Range("G161").Value = 509

$ws.Range("A162:G162").PasteSpecial(-4122)
$ws.Range("A162").Value = 44123
$ws.Range("B162").Value = "10:15:00"
$ws.Range("C162").Value = 3778
$ws.Range("D162").Value = 1838
$ws.Range("E162").Value = 110
$ws.Range("F162").Value = 3016
$ws.Range("G162").Value = 652

$ws.Range("A163:G163").PasteSpecial(-4122)
$ws.Range("A163").Value = 44124
$ws.Range("B163").Value = "09:15:00"
$ws.Range("C163").Value = 3816
$ws.Range("D163").Value = 1855
$ws.Range("E163").Value = 110
$ws.Range("F163").Value = 3084
$ws.Range("G163").Value = 622

$ws.Range("A164:G164").PasteSpecial(-4122)
$ws.Range("A164").Value = 44125
$ws.Range("B164").Value = "10:30:00"
$ws.Range("C164").Value = 3912
$ws.Range("D164").Value = 1894
$ws.Range("E164").Value = 110
$ws.Range("F164").Value = 3157
$ws.Range("G164").Value = 645

$ws.Range("A165:G165").PasteSpecial(-4122)
$ws.Range("A165").Value = 44126
$ws.Range("B165").Value = "09:30:00"
$ws.Range("C165").Value = 4103
$ws.Range("D165").Value = 1973
$ws.Range("E165").Value = 111
$ws.Range("F165").Value = 3237
$ws.Range("G165").Value = 755

$ws.Range("A166:G166").PasteSpecial(-4122)
$ws.Range("A166").Value = 44127
$ws.Range("B166").Value = "10:30:00"
$ws.Range("C166").Value = 4288
$ws.Range("D166").Value = 2031
$ws.Range("E166").Value = 111
$ws.Range("F166").Value = 3330
$ws.Range("G166").Value = 847

$excel.CutCopyMode = $false
